$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "https://fmfb.com.af"
$ws.Range("B4").Value = "info@fmfb.com.af"
$ws.Range("A5").Value = "https://sanayee.org.af"
$ws.Range("B5").Value = "info@sanayee.org.af"
